# U12 and U19 B-side connections to H-bridges defined.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS Dualchip")
# ($wb.ActiveSheet resolves to the same sheet - it's the workbook's active tab.)

# The "IC_*" signal names that were parked in column E (as notes) now
# become the real "X3 Signal" values in column P, replacing the stale
# U12_A#/U19_A# placeholders. The now-redundant column-E notes are cleared.

# P16:P21 -> U19 diag signals
$ws.Range("P16").Value = "IC1_DIAGB"
$ws.Range("P17").Value = "IC1_DIAGA"
$ws.Range("P18").Value = "IC2_DIAGB"
$ws.Range("P19").Value = "IC2_DIAGA"
$ws.Range("P20").Value = "IC3_DIAGB"
$ws.Range("P21").Value = "IC3_DIAGA"

# P28:P33 -> U12 in signals
$ws.Range("P28").Value = "IC1_INB"
$ws.Range("P29").Value = "IC1_INA"
$ws.Range("P30").Value = "IC2_INB"
$ws.Range("P31").Value = "IC2_INA"
$ws.Range("P32").Value = "IC3_INB"
$ws.Range("P33").Value = "IC3_INA"

# Clear the now-redundant column E notes that held these signal names
$ws.Range("E22").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("E34").ClearContents()
$ws.Range("E35").ClearContents()
$ws.Range("E42").ClearContents()
$ws.Range("E43").ClearContents()
$ws.Range("E44").ClearContents()
$ws.Range("E45").ClearContents()

# Update the view's scroll position and selection to match
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("P21").Select()
